$wb = $excel.ActiveWorkbook

# --- Sheet "content": part names placed into destination wells ---
$wsContent = $wb.Worksheets.Item("content")
$wsContent.Range("D2").Value = "p19_mtagbfp2"
$wsContent.Range("E2").Value = "p24_3'ha_haavs1"
$wsContent.Range("F2").Value = "p6_nt-igkl sequence"
$wsContent.Range("G2").Value = "p8_p2a"
$wsContent.Range("D3").Value = "p19_tet-on-3g"
$wsContent.Range("E3").Value = "p24_3'pb"
$wsContent.Range("F3").Value = "p6_nt-mls"
$wsContent.Range("G3").Value = "p8a_ct-kdel"
$wsContent.Range("D4").Value = "p1_5'-itr-pb"
$wsContent.Range("E4").Value = "p25_sv40-ori"
$wsContent.Range("F4").Value = "p6_nt-myristoylation signal"
$wsContent.Range("G4").Value = "p8b_ires2"
$wsContent.Range("D5").Value = "p1_5'ha-haavs1"
$wsContent.Range("E5").Value = "p2_insulatorfb"
$wsContent.Range("F5").Value = "p6_nt-palm sequence"
$wsContent.Range("G5").Value = "p9_a-tubulin"
$wsContent.Range("D6").Value = "p20_ct-minute-nes"
$wsContent.Range("E6").Value = "p3_cagp"
$wsContent.Range("F6").Value = "p6_nt-sv40_nls"
$wsContent.Range("G6").Value = "p9_bsdr"
$wsContent.Range("D7").Value = "p20_ct-nes"
$wsContent.Range("E7").Value = "p3_cmvp_tet"
$wsContent.Range("F7").Value = "p7_bxb1"
$wsContent.Range("G7").Value = "p9_dmrc"
$wsContent.Range("D8").Value = "p20_linker3"
$wsContent.Range("E8").Value = "p3_ef1ap"
$wsContent.Range("F8").Value = "p7_l7ae-weiss"
$wsContent.Range("G8").Value = "p9_firefly_luciferase"
$wsContent.Range("D9").Value = "p20_p2a"
$wsContent.Range("E9").Value = "p3_tre3gp"
$wsContent.Range("F9").Value = "p7_l7ae"
$wsContent.Range("G9").Value = "p9_mneogreen"
$wsContent.Range("D10").Value = "p21_dmra"
$wsContent.Range("E10").Value = "p4_kt-l7ae -weiss"
$wsContent.Range("F10").Value = "p7_laci"
$wsContent.Range("G10").Value = "p9_mruby2"
$wsContent.Range("D11").Value = "p21_mkate2"
$wsContent.Range("E11").Value = "p4_lac-o"
$wsContent.Range("F11").Value = "p7_mcherry"
$wsContent.Range("G11").Value = "p9_mtagbfp2"
$wsContent.Range("D12").Value = "p21_mneogreen"
$wsContent.Range("E12").Value = "p5_attb-bxb1"
$wsContent.Range("F12").Value = "p7_mkate2"
$wsContent.Range("G12").Value = "p9_neor"
$wsContent.Range("D13").Value = "p21_mruby2"
$wsContent.Range("E13").Value = "p5_attp-bxb1"
$wsContent.Range("F13").Value = "p7_mneogreen"
$wsContent.Range("G13").Value = "p9_puror"
$wsContent.Range("D14").Value = "p21_mtagbfp2"
$wsContent.Range("E14").Value = "p5_k1-k1"
$wsContent.Range("F14").Value = "p7_mruby2"
$wsContent.Range("D15").Value = "p21_puror"
$wsContent.Range("E15").Value = "p5_kt-weiss"
$wsContent.Range("F15").Value = "p7_mtagbfp2"
$wsContent.Range("D16").Value = "p22_pgkpolya"
$wsContent.Range("E16").Value = "p6_atg_boxc"
$wsContent.Range("F16").Value = "p8_linker1"
$wsContent.Range("Y16").Value = "WATER"
$wsContent.Range("D17").Value = "p23_insulatorfb"
$wsContent.Range("E17").Value = "p6_kozak-atg"
$wsContent.Range("F17").Value = "p8_linker2"
$wsContent.Range("Y17").Value = "BUFFER"

# --- Sheet "volume (uL)": transfer volumes ---
$wsVolume = $wb.Worksheets.Item("volume (uL)")
$wsVolume.Range("D2").Value = 50
$wsVolume.Range("E2").Value = 50
$wsVolume.Range("F2").Value = 50
$wsVolume.Range("G2").Value = 50
$wsVolume.Range("D3").Value = 50
$wsVolume.Range("E3").Value = 50
$wsVolume.Range("F3").Value = 50
$wsVolume.Range("G3").Value = 50
$wsVolume.Range("D4").Value = 50
$wsVolume.Range("E4").Value = 50
$wsVolume.Range("F4").Value = 50
$wsVolume.Range("G4").Value = 50
$wsVolume.Range("D5").Value = 50
$wsVolume.Range("E5").Value = 50
$wsVolume.Range("F5").Value = 50
$wsVolume.Range("G5").Value = 50
$wsVolume.Range("D6").Value = 50
$wsVolume.Range("E6").Value = 50
$wsVolume.Range("F6").Value = 50
$wsVolume.Range("G6").Value = 50
$wsVolume.Range("D7").Value = 50
$wsVolume.Range("E7").Value = 50
$wsVolume.Range("F7").Value = 50
$wsVolume.Range("G7").Value = 50
$wsVolume.Range("D8").Value = 50
$wsVolume.Range("E8").Value = 50
$wsVolume.Range("F8").Value = 50
$wsVolume.Range("G8").Value = 50
$wsVolume.Range("D9").Value = 50
$wsVolume.Range("E9").Value = 50
$wsVolume.Range("F9").Value = 50
$wsVolume.Range("G9").Value = 50
$wsVolume.Range("D10").Value = 50
$wsVolume.Range("E10").Value = 50
$wsVolume.Range("F10").Value = 50
$wsVolume.Range("G10").Value = 50
$wsVolume.Range("D11").Value = 50
$wsVolume.Range("E11").Value = 50
$wsVolume.Range("F11").Value = 50
$wsVolume.Range("G11").Value = 50
$wsVolume.Range("D12").Value = 50
$wsVolume.Range("E12").Value = 50
$wsVolume.Range("F12").Value = 50
$wsVolume.Range("G12").Value = 50
$wsVolume.Range("D13").Value = 50
$wsVolume.Range("E13").Value = 50
$wsVolume.Range("F13").Value = 50
$wsVolume.Range("G13").Value = 50
$wsVolume.Range("D14").Value = 50
$wsVolume.Range("E14").Value = 50
$wsVolume.Range("F14").Value = 50
$wsVolume.Range("D15").Value = 50
$wsVolume.Range("E15").Value = 50
$wsVolume.Range("F15").Value = 50
$wsVolume.Range("D16").Value = 50
$wsVolume.Range("E16").Value = 50
$wsVolume.Range("F16").Value = 50
$wsVolume.Range("D17").Value = 50
$wsVolume.Range("E17").Value = 50
$wsVolume.Range("F17").Value = 50
$wsVolume.Range("Y16").Value = 50
$wsVolume.Range("Y17").Value = 50

# --- Sheet "concentration (ng-uL)": source concentrations ---
$wsConc = $wb.Worksheets.Item("concentration (ng-uL)")
$wsConc.Range("D2").Value = 68.7
$wsConc.Range("E2").Value = 78.6
$wsConc.Range("F2").Value = 69.6
$wsConc.Range("G2").Value = 68.7
$wsConc.Range("D3").Value = 98.3
$wsConc.Range("E3").Value = 52.8
$wsConc.Range("F3").Value = 86.00000000000001
$wsConc.Range("G3").Value = 98.3
$wsConc.Range("D4").Value = 68.3
$wsConc.Range("E4").Value = 46.4
$wsConc.Range("F4").Value = 73.7
$wsConc.Range("G4").Value = 68.3
$wsConc.Range("D5").Value = 80.3
$wsConc.Range("E5").Value = 59.3
$wsConc.Range("F5").Value = 69.5
$wsConc.Range("G5").Value = 80.3
$wsConc.Range("D6").Value = 111.1
$wsConc.Range("E6").Value = 71.6
$wsConc.Range("F6").Value = 56.1
$wsConc.Range("G6").Value = 111.1
$wsConc.Range("D7").Value = 76.3
$wsConc.Range("E7").Value = 106.7
$wsConc.Range("F7").Value = 63.4
$wsConc.Range("G7").Value = 76.3
$wsConc.Range("D8").Value = 77.3
$wsConc.Range("E8").Value = 70.3
$wsConc.Range("F8").Value = 44.6
$wsConc.Range("G8").Value = 77.3
$wsConc.Range("D9").Value = 80.1
$wsConc.Range("E9").Value = 68.7
$wsConc.Range("F9").Value = 62.40000000000001
$wsConc.Range("G9").Value = 80.1
$wsConc.Range("D10").Value = 50.1
$wsConc.Range("E10").Value = 98.3
$wsConc.Range("F10").Value = 226
$wsConc.Range("G10").Value = 50.1
$wsConc.Range("D11").Value = 78.6
$wsConc.Range("E11").Value = 68.3
$wsConc.Range("F11").Value = 78.6
$wsConc.Range("G11").Value = 78.6
$wsConc.Range("D12").Value = 52.8
$wsConc.Range("E12").Value = 80.3
$wsConc.Range("F12").Value = 52.8
$wsConc.Range("G12").Value = 52.8
$wsConc.Range("D13").Value = 46.4
$wsConc.Range("E13").Value = 111.1
$wsConc.Range("F13").Value = 46.4
$wsConc.Range("G13").Value = 46.4
$wsConc.Range("D14").Value = 59.3
$wsConc.Range("E14").Value = 76.3
$wsConc.Range("F14").Value = 59.3
$wsConc.Range("D15").Value = 71.6
$wsConc.Range("E15").Value = 77.3
$wsConc.Range("F15").Value = 71.6
$wsConc.Range("D16").Value = 106.7
$wsConc.Range("E16").Value = 80.1
$wsConc.Range("F16").Value = 106.7
$wsConc.Range("Y16").Value = 1
$wsConc.Range("D17").Value = 70.3
$wsConc.Range("E17").Value = 50.1
$wsConc.Range("F17").Value = 70.3
$wsConc.Range("Y17").Value = 1
